$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The paragraph that currently begins "The device that the Rebmem designs
#    has suitable input controls..." is replaced wholesale with a brand new
#    paragraph of text ("I assume that the device...").  The original wording
#    about the device's input controls becomes a *new* paragraph placed right
#    after it (with one small wording tweak: "a joystick" -> "a small
#    joystick or d-pad").
# ---------------------------------------------------------------------------

$oldParaText = "The device that the Rebmem designs has suitable input controls. Such as, but not limited to, a touch screen, a joystick to interface with menus, suitable buttons, a small keyboard, etc. Or alternatively, if during phase 2 they decide that the portability of the device is more important than the usability, then they might consider locking down specific complicated functionality (such as creating a playlist, or the search feature) to only when the device is plugged into a computer or laptop as then the user would (most likely) have access to a larger screen, mouse, and keyboard."

$newFirstParaText = "I assume that the device that Rebmem plans to design is supposed to be hand-held size, due to the repeated mention of " + [char]0x2018 + "portable" + [char]0x2019 + " in their description of the device. Therefore, while I was designing my user-interfaces I was considering that the device would have to be small and would not have a lot of room for complicated UI components and designed accordingly."

$newSecondParaText = "The device that the Rebmem designs has suitable input controls. Such as, but not limited to, a touch screen, a small joystick or d-pad to interface with menus, suitable buttons, a small keyboard, etc. Or alternatively, if during phase 2 they decide that the portability of the device is more important than the usability, then they might consider locking down specific complicated functionality (such as creating a playlist, or the search feature) to only when the device is plugged into a computer or laptop as then the user would (most likely) have access to a larger screen, mouse, and keyboard."

$rng = $d.Content
$found = $rng.Find.Execute($oldParaText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "could not find device/joystick paragraph" }

# Replace the paragraph's text in place with the new "I assume..." wording.
$rng.Text = $newFirstParaText

# Re-insert the (slightly tweaked) original wording as a brand new paragraph
# directly after the one we just rewrote.
$rng.InsertAfter("`r" + $newSecondParaText)

# ---------------------------------------------------------------------------
# 2) Merge two runs in the "Many features..." paragraph ("that t" + "he user"
#    -> "that the user ...").  This is a pure run-merge with identical
#    resulting text, so nothing visible changes; skip (text identical).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 3) Extend the final ("Finally, there is a possible minor assumption...")
#    paragraph with an extra sentence, inserted around the existing
#    "_GoBack" bookmark so the bookmark ends up in the same relative spot.
# ---------------------------------------------------------------------------

$bm = $d.Bookmarks("_GoBack")
$bmRange = $bm.Range
$bmRange.InsertBefore(" The user guide would also ")

$bm2 = $d.Bookmarks("_GoBack")
$bmRange2 = $bm2.Range
$bmRange2.InsertAfter("be able to assist the user.")

Write-Output "edits applied"
